$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: average of the J column (|S*|/n)
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"

# Summary rows 14-17
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Build the bold/size-12/vertically-centered style once on a scratch cell,
# then copy its format onto B14:B17 in a single operation so we don't leave
# unused intermediate styles behind in the style table.
$tmp = $ws.Range("AA1")
$tmp.Font.Bold = $true
$tmp.Font.Size = 12
$tmp.VerticalAlignment = -4108

$rng = $ws.Range("B14:B17")
$tmp.Copy()
$rng.PasteSpecial(-4122)
$tmp.Clear()

# Row heights for the new summary rows
$ws.Range("A14:B17").RowHeight = 15.6

# Match the selection left behind in the saved file
[void]$ws.Range("A14:B17").Select()

# Page setup (paper size 9 = A4, portrait orientation)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
